# "Generate Report for Handback"
#
# The localization-status report has two data rows per language sheet
# (one per source file). Row 2 corresponds to
# "1df68d2c-6b5c-4aad-92ca-a6905a981c51.md". A new handback was generated
# for that file, so its "Latest Handback DateTime" column (K) is refreshed
# on both the "zh-cn" and "de-de" sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("K2").Value = "2016-10-14 08:31:37"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("K2").Value = "2016-10-14 08:31:53"
